$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing the existing header/data/totals down
# by one row (matches the row-shift seen throughout the diff).
$ws.Rows.Item(1).Insert()

# Move the "done" markers from column C to column B for the left-hand
# ("ongoing") block (rows 3-11), and populate the new right-hand
# ("done") block's marker column H.
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 8).Value = 1
}

# New labelled, merged, centre-aligned header cells on row 1.
$ws.Range("B1").Value = "ongoing"
$ws.Range("G1").Value = "done"

$ws.Range("B1:D1").Merge()
$ws.Range("G1:I1").Merge()

$ws.Range("B1:D1").HorizontalAlignment = -4108
$ws.Range("G1:I1").HorizontalAlignment = -4108

# Match the selection shown in the edited file.
$ws.Range("G2").Select() | Out-Null

$wb.Save() | Out-Null
